# T3.Calc.xlsx edit: rename the "index" column to "i" and convert it
# from a 1-based row index to a 0-based row index; adjust the column A
# width; update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "index" -> "i" (this also renames the "testdata" table's
# first column, since it is the header cell of that table).
$ws.Range("A1").Value2 = "i"

# Re-number the index column from 1-based (1..502) to 0-based (0..501)
# for data rows 2..503.
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value2
    $cell.Value2 = $cur - 1
}

# Narrow column A now that the values have fewer digits.
# (ColumnWidth uses Excel's internal character-width units; 3.14 rounds
# to the same stored "width" of 4 that Excel itself produces here.)
$ws.Columns.Item(1).ColumnWidth = 3.14

# Update the active selection to L12, as recorded in the saved view state.
$ws.Range("L12").Select()
